# Comments on Tim's models
#
# Adds three of "Mike's Comments" review notes to the document:
#   1. In the paragraph that holds the first tab stop (pos 7830) / tab run
#      right after the intro picture, insert a new run of text before the
#      existing tab run.
#   2. The first otherwise-empty paragraph that follows that same tab-stop
#      group gets a new run with a comment about the diagram titles.
#   3. At the very end of the document (before the section break), append
#      eight new paragraphs that match the formatting (tab stop pos 2550)
#      of the document's final paragraph; the last of those new paragraphs
#      carries a comment about the incremental-design picture.

$d = $word.ActiveDocument
$rightQuote = [char]0x2019

function Get-CustomTabPositions($paragraph) {
    $positions = @()
    $tabStops = $paragraph.Format.TabStops
    for ($j = 1; $j -le $tabStops.Count; $j++) {
        $stop = $tabStops.Item($j)
        if ($stop.CustomTab) {
            $positions += $stop.Position
        }
    }
    return $positions
}

$count = $d.Paragraphs.Count

# ---------------------------------------------------------------------
# Change 1 + 2: locate the run of paragraphs that share the 391.5pt
# (7830 twip) custom tab stop. The first paragraph of that run gets a
# new leading comment run (before its existing tab run); the paragraph
# immediately after the run ends gets a brand-new comment run of its
# own (it currently has no runs at all).
# ---------------------------------------------------------------------
$firstTabbed = -1
$lastTabbed = -1
for ($i = 1; $i -le $count; $i++) {
    $para = $d.Paragraphs.Item($i)
    $positions = Get-CustomTabPositions $para
    if ($positions -contains 391.5) {
        if ($firstTabbed -eq -1) {
            $firstTabbed = $i
        }
        $lastTabbed = $i
    }
}

if ($firstTabbed -ne -1) {
    $introPara = $d.Paragraphs.Item($firstTabbed)
    $insertionPoint = $d.Range($introPara.Range.Start, $introPara.Range.Start)
    $insertionPoint.InsertBefore("Mike" + $rightQuote + "s Comments: Add bullet point list of requirements and specific details of each step")
}

if ($lastTabbed -ne -1) {
    $diagramPara = $d.Paragraphs.Item($lastTabbed + 1)
    $diagramPara.Range.Text = "Mike" + $rightQuote + "s Comments: Add titles of what each diagram represents. What type of model is this representing?"
}

# ---------------------------------------------------------------------
# Change 3: append eight new paragraphs (inheriting the final
# paragraph's formatting) at the end of the document; the last one
# carries the comment about the incremental-design picture.
# ---------------------------------------------------------------------
for ($k = 0; $k -lt 8; $k++) {
    $n = $d.Paragraphs.Count
    $lastPara = $d.Paragraphs.Item($n)
    $lastPara.Range.InsertParagraphAfter()
}

$finalCount = $d.Paragraphs.Count
$finalPara = $d.Paragraphs.Item($finalCount)
$finalPara.Range.Text = "Mike" + $rightQuote + "s Comments: What is the basic product look like and what changes will be made with each incremental design?"
